# Roles.xlsx: insert a new "scripts" column before the existing "ID" column.
# Existing column F ("ID") and G ("mod") shift right to G and H respectively,
# and the new column F is populated with per-role script strings.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column at F, shifting ID/mod (old F/G) to G/H.
$ws.Columns("F:F").Insert()

# Header
$ws.Range("F1").Value = "scripts"

# Per-row "scripts" values (row 3 / Mage intentionally left blank).
$ws.Range("F2").Value = "prod_speed,25"
$ws.Range("F4").Value = "fertility,20 | recipe_speed,50,maid_recipe"
$ws.Range("F5").Value = "manipulation,20 | prod_speed,25"
$ws.Range("F6").Value = "sight,20 | prod_speed,25"
$ws.Range("F7").Value = "compliance,20 | prod_speed,25"
$ws.Range("F8").Value = "manipulation,20 | prod_speed,25"
$ws.Range("F9").Value = "speech,20 | prod_speed,25"
$ws.Range("F10").Value = "prod_speed,25"
$ws.Range("F11").Value = "fertility,20 | recipe_speed,50,tavern_recipe"
$ws.Range("F12").Value = "focus,20 | prod_speed,25"
$ws.Range("F13").Value = "recipe_speed,50,dryad_recipe"
$ws.Range("F14").Value = "speech,20 | recipe_speed,50,overseer_recipe"
$ws.Range("F15").Value = "manipulation,20 | building_prod_speed,50,mine"
$ws.Range("F16").Value = "fertility,20 | prod_speed,25"
$ws.Range("F17").Value = "focus,20 | prod_speed,25"
$ws.Range("F18").Value = "manipulation,20 | prod_speed,25"
$ws.Range("F19").Value = "fertility,20 | recipe_speed,50,cow_recipe"
$ws.Range("F20").Value = "moving,20 | prod_speed,25"
$ws.Range("F21").Value = "manipulation,20 | prod_speed,25"
$ws.Range("F22").Value = "prod_speed,25"
$ws.Range("F23").Value = "manipulation,20 | prod_speed,25"
$ws.Range("F24").Value = "move_speed,20 | initiative,1 | prod_speed,25"
$ws.Range("F25").Value = "moving,20 | prod_speed,25"
$ws.Range("F26").Value = "focus,20 | recipe_speed,50,wool_recipe"
$ws.Range("F27").Value = "IF:combat_terrain,snow | damage_bonus,1 | prod_speed,25"
$ws.Range("F28").Value = "moving,20 | prod_speed,25"
$ws.Range("F29").Value = "IF:combat_terrain,aphrodisiac | damage_bonus,1 | prod_speed,25"
$ws.Range("F30").Value = "IF:combat_terrain,latex_pool | damage_bonus,1 | prod_speed,25"
$ws.Range("F31").Value = "fertility,20 | prod_speed,25"
$ws.Range("F32").Value = "recipe_speed,50,metal_slime_recipe"
$ws.Range("F33").Value = "focus,20 | prod_speed,25"
$ws.Range("F34").Value = "focus,25 | recipe_speed,50,silk_recipe"
$ws.Range("F35").Value = "focus,20 | prod_speed,25"
$ws.Range("F36").Value = "moving,20 | prod_speed,25"
$ws.Range("F37").Value = "speech,20 | prod_speed,25"
$ws.Range("F38").Value = "manipulation,10 | prod_speed,25"
$ws.Range("F39").Value = "IF:combat_terrain,snow | damage_reduction,1 | prod_speed,25"
$ws.Range("F40").Value = "prod_speed,25 | prod_speed,25"
$ws.Range("F41").Value = "move_speed,20 | prod_speed,25"
$ws.Range("F42").Value = "focus,10 | prod_speed,25"
$ws.Range("F43").Value = "move_speed,20 | prod_speed,25"
$ws.Range("F44").Value = "fertility,10 | prod_speed,25"
$ws.Range("F45").Value = "move_speed,20 | prod_speed,25"
$ws.Range("F46").Value = "manipulation,10 | prod_speed,25"
$ws.Range("F47").Value = "prod_speed,25"
$ws.Range("F48").Value = "manipulation,10 | prod_speed,25"
$ws.Range("F49").Value = "manipulation,10 | prod_speed,25"

Write-Host "scripts column inserted and populated"
